$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a numeric-looking value as TEXT (matches the source file's
# convention where every data cell is stored as a string, even when the
# content happens to look like a number). We briefly force a Text number
# format so Excel's auto-type-detection doesn't coerce the entry into a
# real number, then strip the format back off so the cell doesn't end up
# with a residual style index.
function Set-TextValue($addr, $value) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.ClearFormats()
}

# Helper: write a literal empty string into a cell (as opposed to simply
# clearing it, which would drop the cell from the sheet entirely). A lone
# leading apostrophe forces Excel to commit an empty-string text cell; we
# then drop the resulting quote-prefix style the same way as above.
function Set-EmptyStringValue($addr) {
    $rng = $ws.Range($addr)
    $rng.Value = "'"
    $rng.ClearFormats()
}

# Row 29
Set-TextValue "C29" "53490059"
$ws.Range("D29").Value = "coluna  do conj transversal traseiro ld"
Set-TextValue "E29" "21"
$ws.Range("F29").Value = "FERRAMENTARIA"
$ws.Range("G29").Value = "matheus"
$ws.Range("H29").Value = "2º TURNO"
$ws.Range("I29").Value = "CMM GLOBAL"
$ws.Range("J29").Value = "LAMENTAÇÃO CLIENTE"
Set-EmptyStringValue "K29"
$ws.Range("L29").Value = "C2025.0027"

# Row 30
Set-TextValue "C30" "53437117"
$ws.Range("D30").Value = "SOLITÁRIA LD"
Set-TextValue "E30" "21"
$ws.Range("F30").Value = "ENGENHARIA"
$ws.Range("G30").Value = "matheus"
$ws.Range("H30").Value = "2º TURNO"
$ws.Range("I30").Value = "PAQUÍMETRO"
$ws.Range("J30").Value = "ODM"
Set-EmptyStringValue "K30"
$ws.Range("L30").Value = "C2025.0028"
